# The workbook previously had a "lib" of leftover Apache POI-generated
# shared-string entries ("ParentEmail" / "RequestVolunteer") duplicated
# without spaces. This cleans them up by renaming the header cells on the
# "Students" sheet to their properly spaced versions, which collapses the
# shared-string table (removing the now-unused duplicate entries) and
# updates the active cell selection on both data sheets.

$wb = $excel.ActiveWorkbook

$wsStudents = $wb.Worksheets.Item("Students")
$wsVolunteers = $wb.Worksheets.Item("Volunteers")

# Rename header cells to the corrected, spaced labels.
$wsStudents.Range("G1").Value = "Parent Email"
$wsStudents.Range("K1").Value = "Requested Volunteer"

# Update the selected cell on each sheet to match the saved view state.
$wsStudents.Range("L1").Select()
$wsVolunteers.Range("K1").Select()
